# Apply updates described by the commit diff to the "Artfynd" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 62: a handful of discrete field edits (record id, sex, activity, and
# the public comment text).
# ---------------------------------------------------------------------------
$ws.Range("A62").Value2 = 80976130
$ws.Range("L62").Value2 = "hane"
$ws.Range("M62").Value2 = "frispringande/krypande"
$ws.Range("AC62").Value2 = "Ny lokal, och blott tredje kända i Södermanland!? Närmast funnen i Tyresta NP. Bör eftersökas på fler lokaler i kommunen! Grävde först fram fragment av en ad hona, men kunde sedan finna en vuxen hane i en perfekt rödmurken granlåga i sent nedbrytningsstadium."

# ---------------------------------------------------------------------------
# Rows 64-66: the three records got cyclically re-ordered (row 64 now holds
# what used to be row 65's data, row 65 now holds what used to be row 66's
# data, and row 66 now holds what used to be row 64's data). Capture each
# full row first, then write them back rotated.
#
# Some columns store values that look numeric but are meant to stay text
# (I/J hold counts such as "3"/"fruktkroppar", and Y/Z/AA/AB hold
# dates/times such as "2023-08-25"/"09:36"). Force those columns to a text
# number format before writing so Excel doesn't reinterpret the strings as
# numbers or date/time serial values.
# ---------------------------------------------------------------------------
$ws.Range("I64:J66").NumberFormat = "@"
$ws.Range("Y64:AB66").NumberFormat = "@"

$row64 = $ws.Range("A64:AY64").Value2
$row65 = $ws.Range("A65:AY65").Value2
$row66 = $ws.Range("A66:AY66").Value2

$ws.Range("A64:AY64").Value2 = $row65
$ws.Range("A65:AY65").Value2 = $row66
$ws.Range("A66:AY66").Value2 = $row64
